$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 28
$ws.Range("C3").Value = 40
$ws.Range("C4").Value = 27
$ws.Range("C5").Value = 24
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 102
$ws.Range("C8").Value = 71
$ws.Range("C9").Value = 23
$ws.Range("C10").Value = 56
$ws.Range("C11").Value = 12
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 25
$ws.Range("C14").Value = 26
$ws.Range("C15").Value = 74
$ws.Range("C16").Value = 175
$ws.Range("C17").Value = 14
$ws.Range("C18").Value = 7
$ws.Range("C19").Value = 149
$ws.Range("C20").Value = 59
$ws.Range("C21").Value = 115
$ws.Range("C23").Value = 125
$ws.Range("C24").Value = 69
$ws.Range("C25").Value = 83
$ws.Range("C26").Value = 57
$ws.Range("C27").Value = 43
$ws.Range("C28").Value = 68
$ws.Range("C29").Value = 91
$ws.Range("C30").Value = 82
$ws.Range("C31").Value = 15
$ws.Range("C32").Value = 42
$ws.Range("C33").Value = 32
$ws.Range("C34").Value = 61
$ws.Range("C35").Value = 30
$ws.Range("C36").Value = 123
$ws.Range("C37").Value = 81
$ws.Range("C38").Value = 39
$ws.Range("C39").Value = 45
$ws.Range("C40").Value = 19
$ws.Range("C41").Value = 18
$ws.Range("C42").Value = 126
